$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GVA")

# Row 2 (Qminus1)
$ws.Range("B2").Value = 0.006109567252182264
$ws.Range("C2").Value = 0.4815469527244732
$ws.Range("D2").Value = 0.3941471084727503
$ws.Range("E2").Value = 0.6278113637652238
$ws.Range("F2").Value = 0.6514797914836403

# Row 3 (Q0)
$ws.Range("B3").Value = 0.04939103104625085
$ws.Range("C3").Value = 0.5032213284952747
$ws.Range("D3").Value = 0.3769116481242295
$ws.Range("E3").Value = 0.6139313057046607
$ws.Range("F3").Value = 0.633419439721888
$ws.Range("G3").Value = 15

# Row 4 (Q1)
$ws.Range("B4").Value = 0.3050745724485351
$ws.Range("C4").Value = 0.781138290135054
$ws.Range("D4").Value = 0.7388866705223734
$ws.Range("E4").Value = 0.8595851735124178
$ws.Range("F4").Value = 0.8339630346331093
$ws.Range("G4").Value = 14
